# Adds three new observation rows (4, 5, 6) to the "Artfynd" sheet,
# matching the Artportalen-style export rows already present in rows 2-3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make sure date/time-looking text columns (Y, Z, AA, AB) are stored as
# plain text, not auto-converted to date/time serial numbers.
$ws.Range("Y4:AB6").NumberFormat = "@"

function Set-Row($r, $id, $taxonSortOrder, $valideringsstatus, $rodlistade, $taxonId, $artnamn, $vetenskapligtNamn, $auktor, $lokalnamn, $ost, $nord, $noggrannhet, $lan, $kommun, $provins, $forsamling, $startdatum, $starttid, $slutdatum, $sluttid, $rapportor) {
    $ws.Cells.Item($r, 1).Value = $id
    $ws.Cells.Item($r, 2).Value = $taxonSortOrder
    $ws.Cells.Item($r, 3).Value = $valideringsstatus
    $ws.Cells.Item($r, 4).Value = $rodlistade
    $ws.Cells.Item($r, 5).Value = $taxonId
    $ws.Cells.Item($r, 6).Value = $artnamn
    $ws.Cells.Item($r, 7).Value = $vetenskapligtNamn
    $ws.Cells.Item($r, 8).Value = $auktor
    $ws.Cells.Item($r, 16).Value = $lokalnamn
    $ws.Cells.Item($r, 17).Value = $ost
    $ws.Cells.Item($r, 18).Value = $nord
    $ws.Cells.Item($r, 19).Value = $noggrannhet
    $ws.Cells.Item($r, 20).Value = $lan
    $ws.Cells.Item($r, 21).Value = $kommun
    $ws.Cells.Item($r, 22).Value = $provins
    $ws.Cells.Item($r, 23).Value = $forsamling
    $ws.Cells.Item($r, 25).Value = $startdatum
    $ws.Cells.Item($r, 26).Value = $starttid
    $ws.Cells.Item($r, 27).Value = $slutdatum
    $ws.Cells.Item($r, 28).Value = $sluttid
    $ws.Cells.Item($r, 30).Value = $False
    $ws.Cells.Item($r, 31).Value = $False
    $ws.Cells.Item($r, 33).Value = $False
    $ws.Cells.Item($r, 49).Value = $rapportor
    $ws.Cells.Item($r, 50).Value = $rapportor
}

Set-Row 4 112128524 90666 "Ovaliderad" `
    "LC" 4364 "Dropptaggsvamp" `
    "Hydnellum ferrugineum" "(Fr.:Fr.) P. Karst." `
    "Godmyr (Godmyr), Ly lm" 690280.6233055658 7126404.088587272 `
    1 "Västerbotten" "Lycksele" "Lycksele lappmark" `
    "Örträsk" "2023-09-16" "15:00" `
    "2023-09-16" "15:00" "Ulrika Karlsson"

Set-Row 5 112127587 90660 "Ovaliderad" `
    "NT" 4362 "Blå taggsvamp" `
    "Hydnellum caeruleum" "(Hornem.) P.Karst." `
    "Svarvarmyran (Svarvarmyran), Ly lm" 690447.4141232002 7125628.52413491 `
    1 "Västerbotten" "Lycksele" "Lycksele lappmark" `
    "Örträsk" "2023-09-16" "14:27" `
    "2023-09-16" "14:27" "Ulrika Karlsson"

Set-Row 6 112127546 90658 "Ovaliderad" `
    "NT" 4361 "Orange taggsvamp" `
    "Hydnellum aurantiacum" "(Batsch:Fr.) P.Karst." `
    "Svarvarmyran (Svarvarmyran), Ly lm" 690408.1690133474 7125570.209357738 `
    1 "Västerbotten" "Lycksele" "Lycksele lappmark" `
    "Örträsk" "2023-09-16" "14:25" `
    "2023-09-16" "14:25" "Ulrika Karlsson"
